$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("new")

# Update the four cells in row 4 with refined values
$ws.Range("H4").Value = 1.03
$ws.Range("I4").Value = 2.63
$ws.Range("J4").Value = 0.49
$ws.Range("K4").Value = 2.14

# Move the active selection to K5 on this sheet
$ws.Activate()
$ws.Range("K5").Select()
